$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "ADs" sheet: the four most recent entries (rows 5:8 - Daeloan, kafka,
#    Brilith, Mr.Skieyagi69) have already been folded into the new "Cities"
#    print-out below, so they are removed from here to keep the two outputs
#    in sync.
# ---------------------------------------------------------------------------
$adsSheet = $wb.Worksheets.Item("ADs")
$adsSheet.Rows("5:8").Delete()

# ---------------------------------------------------------------------------
# 2. "Cities" sheet: add the fresh Yamagata "dm" entries that mirror what was
#    just printed out, each with its ID lookup and day-count formula.
# ---------------------------------------------------------------------------
$citiesSheet = $wb.Worksheets.Item("Cities")

$rows = @(
    @{ Row = 2;  Name = "Nyarlathotep"; Pts = 614 },
    @{ Row = 3;  Name = "Alecks_";      Pts = 334 },
    @{ Row = 4;  Name = "Noah";         Pts = 311 },
    @{ Row = 5;  Name = "Aileen";       Pts = 254 },
    @{ Row = 6;  Name = "Letsi";        Pts = 240 },
    @{ Row = 7;  Name = "Nefi85";       Pts = 70 },
    @{ Row = 8;  Name = "xljhx31";      Pts = 58 },
    @{ Row = 9;  Name = "mido009";      Pts = 55 },
    @{ Row = 10; Name = "Hallen98";     Pts = 44 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $citiesSheet.Range("A$row").Formula = "=VLOOKUP(B$row,IDs!B:C,2,FALSE)"
    $citiesSheet.Range("B$row").Value = $r.Name
    $citiesSheet.Range("C$row").Value = $r.Pts
    $citiesSheet.Range("D$row").Value = 2
    $citiesSheet.Range("E$row").Value = "Yamagata"
    $citiesSheet.Range("F$row").Value = "dm"
    $citiesSheet.Range("G$row").Value = "dm"

    $citiesSheet.Range("H$row").NumberFormat = "dd.mm.YYYY"
    $citiesSheet.Range("H$row").Value = 45584.0

    $citiesSheet.Range("I$row").Formula = "=TODAY()-G$row"
}
